$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above the current row 13 ("Batch system"), shifting rows
# 13-14 down to 14-15, matching the diff's row insertion.
$ws.Rows.Item(13).Insert()

# New cell content for the inserted row.
$ws.Range("A13").Value = "MPI executable prefix"

# Copy formatting (font/border/number format) from the row above so the new
# row reuses the existing "label"/"value" cell styles instead of picking up
# blank defaults.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)   # xlPasteFormats

# Match the row height used by its sibling rows (18.5pt).
$ws.Rows.Item(13).RowHeight = 18.5

# Switch the active sheet back to "MAIN Config." and select B13, matching
# the new selection/tabSelected state recorded in the workbook.
$ws.Activate()
$ws.Range("B13").Select()
